$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44956
$ws.Range("M2").Value = 200
$ws.Range("N2").Value = 2000
$ws.Range("O2").Value = 2200
$ws.Range("P2").Value = 2100
$ws.Range("Q2").Value = '$/bandeja 2 kilos'
$ws.Range("R2").Value = 'Provincia de Curicó'
$ws.Range("S2").Value = 1050
$ws.Range("T2").Value = 2

# Row 3
$ws.Range("D3").Value = 44965
$ws.Range("M3").Value = 60
$ws.Range("N3").Value = 2000
$ws.Range("O3").Value = 2200
$ws.Range("P3").Value = 2100
$ws.Range("Q3").Value = '$/bandeja 2 kilos'
$ws.Range("R3").Value = 'Provincia de Curicó'
$ws.Range("S3").Value = 1050
$ws.Range("T3").Value = 2

# Row 4
$ws.Range("D4").Value = 44547
$ws.Range("Q4").Value = '$/bandeja 12 canastillos 125 gramos'
$ws.Range("R4").Value = 'Región del Maule'
$ws.Range("S4").Value = 3500
$ws.Range("T4").Value = 1.5

# Row 5
$ws.Range("D5").Value = 44907
$ws.Range("M5").Value = 400
$ws.Range("N5").Value = 4500
$ws.Range("O5").Value = 5000
$ws.Range("P5").Value = 4750
$ws.Range("S5").Value = 2375

# Row 6
$ws.Range("D6").Value = 44908
$ws.Range("M6").Value = 400
$ws.Range("N6").Value = 4000
$ws.Range("O6").Value = 4400
$ws.Range("P6").Value = 4200
$ws.Range("Q6").Value = '$/bandeja 2 kilos'
$ws.Range("R6").Value = 'Provincia de Curicó'
$ws.Range("S6").Value = 2100
$ws.Range("T6").Value = 2

# Row 7
$ws.Range("D7").Value = 44984
$ws.Range("N7").Value = 2000
$ws.Range("O7").Value = 2200
$ws.Range("P7").Value = 2100
$ws.Range("S7").Value = 1050

# Row 8
$ws.Range("D8").Value = 44981
$ws.Range("M8").Value = 200

# Row 9
$ws.Range("D9").Value = 44925
$ws.Range("M9").Value = 400
$ws.Range("N9").Value = 3000
$ws.Range("O9").Value = 3500
$ws.Range("P9").Value = 3250
$ws.Range("S9").Value = 1625

# Row 10
$ws.Range("D10").Value = 44890
$ws.Range("N10").Value = 8500
$ws.Range("O10").Value = 9000
$ws.Range("P10").Value = 8750
$ws.Range("S10").Value = 4375

# Row 11
$ws.Range("D11").Value = 44994
$ws.Range("M11").Value = 400
$ws.Range("N11").Value = 2300
$ws.Range("O11").Value = 2500
$ws.Range("P11").Value = 2400
$ws.Range("Q11").Value = '$/bandeja 2 kilos'
$ws.Range("R11").Value = 'Provincia de Curicó'
$ws.Range("S11").Value = 1200
$ws.Range("T11").Value = 2

# Row 12
$ws.Range("D12").Value = 44971
$ws.Range("M12").Value = 300
$ws.Range("N12").Value = 2000
$ws.Range("O12").Value = 2200
$ws.Range("P12").Value = 2100
$ws.Range("S12").Value = 1050

# Row 13
$ws.Range("D13").Value = 44918
$ws.Range("M13").Value = 600
$ws.Range("N13").Value = 3000
$ws.Range("O13").Value = 3500
$ws.Range("P13").Value = 3250
$ws.Range("S13").Value = 1625

# Row 14
$ws.Range("D14").Value = 44537
$ws.Range("M14").Value = 400
$ws.Range("N14").Value = 5000
$ws.Range("O14").Value = 5500
$ws.Range("P14").Value = 5250
$ws.Range("Q14").Value = '$/bandeja 12 canastillos 125 gramos'
$ws.Range("R14").Value = 'Región del Maule'
$ws.Range("S14").Value = 3500
$ws.Range("T14").Value = 1.5

# Row 15
$ws.Range("D15").Value = 44974

# Row 16
$ws.Range("D16").Value = 44998
$ws.Range("M16").Value = 120
$ws.Range("N16").Value = 2300
$ws.Range("O16").Value = 2500
$ws.Range("P16").Value = 2400
$ws.Range("Q16").Value = '$/bandeja 2 kilos'
$ws.Range("R16").Value = 'Provincia de Curicó'
$ws.Range("S16").Value = 1200
$ws.Range("T16").Value = 2

# Row 17
$ws.Range("D17").Value = 44876
$ws.Range("N17").Value = 7500
$ws.Range("O17").Value = 8000
$ws.Range("P17").Value = 7750
$ws.Range("Q17").Value = '$/bandeja 12 canastillos 125 gramos'
$ws.Range("S17").Value = 5167
$ws.Range("T17").Value = 1.5

# Row 18
$ws.Range("D18").Value = 44946
$ws.Range("M18").Value = 240
$ws.Range("N18").Value = 2000
$ws.Range("O18").Value = 2200
$ws.Range("P18").Value = 2100
$ws.Range("S18").Value = 1050

# Row 19
$ws.Range("D19").Value = 44897
$ws.Range("M19").Value = 400
$ws.Range("Q19").Value = '$/bandeja 2 kilos'
$ws.Range("S19").Value = 3125
$ws.Range("T19").Value = 2

# Row 20
$ws.Range("D20").Value = 44964

# Row 21
$ws.Range("D21").Value = 44875
$ws.Range("M21").Value = 400
$ws.Range("N21").Value = 7500
$ws.Range("O21").Value = 8000
$ws.Range("P21").Value = 7750
$ws.Range("Q21").Value = '$/bandeja 12 canastillos 125 gramos'
$ws.Range("S21").Value = 5167
$ws.Range("T21").Value = 1.5

# Row 22
$ws.Range("D22").Value = 44159
$ws.Range("L22").Value = 'Segunda'
$ws.Range("M22").Value = 200
$ws.Range("N22").Value = 6500
$ws.Range("O22").Value = 7000
$ws.Range("P22").Value = 6750
$ws.Range("S22").Value = 4500

# Row 23
$ws.Range("D23").Value = 44966
$ws.Range("M23").Value = 100
$ws.Range("N23").Value = 2000
$ws.Range("O23").Value = 2200
$ws.Range("P23").Value = 2100
$ws.Range("S23").Value = 1050

# Row 24
$ws.Range("D24").Value = 44943
$ws.Range("M24").Value = 200
$ws.Range("N24").Value = 2000
$ws.Range("O24").Value = 2200
$ws.Range("P24").Value = 2100
$ws.Range("Q24").Value = '$/bandeja 2 kilos'
$ws.Range("R24").Value = 'Provincia de Colchagua'
$ws.Range("S24").Value = 1050
$ws.Range("T24").Value = 2

# Row 25
$ws.Range("D25").Value = 44895
$ws.Range("M25").Value = 120
$ws.Range("N25").Value = 8000
$ws.Range("O25").Value = 8500
$ws.Range("P25").Value = 8250
$ws.Range("S25").Value = 4125

# Row 26
$ws.Range("D26").Value = 44895
$ws.Range("M26").Value = 300
$ws.Range("N26").Value = 8000
$ws.Range("O26").Value = 8500
$ws.Range("P26").Value = 8250
$ws.Range("S26").Value = 4125

# Row 27
$ws.Range("D27").Value = 44973
$ws.Range("M27").Value = 200

# Row 28
$ws.Range("D28").Value = 44169
$ws.Range("M28").Value = 400
$ws.Range("N28").Value = 5500
$ws.Range("O28").Value = 6000
$ws.Range("P28").Value = 5750
$ws.Range("Q28").Value = '$/bandeja 12 canastillos 125 gramos'
$ws.Range("S28").Value = 3833
$ws.Range("T28").Value = 1.5

# Row 29
$ws.Range("D29").Value = 44533
$ws.Range("N29").Value = 3500
$ws.Range("O29").Value = 3600
$ws.Range("P29").Value = 3550
$ws.Range("Q29").Value = '$/kilo'
$ws.Range("R29").Value = 'Región del Maule'
$ws.Range("S29").Value = 3550
$ws.Range("T29").Value = 1

# Row 30
$ws.Range("D30").Value = 44980
$ws.Range("M30").Value = 200
$ws.Range("N30").Value = 2000
$ws.Range("O30").Value = 2200
$ws.Range("P30").Value = 2100
$ws.Range("S30").Value = 1050

# Row 31
$ws.Range("D31").Value = 44523
$ws.Range("M31").Value = 300
$ws.Range("N31").Value = 3700
$ws.Range("O31").Value = 3800
$ws.Range("P31").Value = 3750
$ws.Range("Q31").Value = '$/kilo'
$ws.Range("R31").Value = 'Región del Maule'
$ws.Range("S31").Value = 3750
$ws.Range("T31").Value = 1

# Row 32
$ws.Range("D32").Value = 44938
$ws.Range("M32").Value = 600

# Row 33
$ws.Range("D33").Value = 44553
$ws.Range("M33").Value = 400
$ws.Range("N33").Value = 5000
$ws.Range("O33").Value = 5500
$ws.Range("P33").Value = 5250
$ws.Range("Q33").Value = '$/bandeja 12 canastillos 125 gramos'
$ws.Range("R33").Value = 'Región del Maule'
$ws.Range("S33").Value = 3500
$ws.Range("T33").Value = 1.5

# Row 36
$ws.Range("D36").Value = 44530
$ws.Range("M36").Value = 160
$ws.Range("N36").Value = 3600
$ws.Range("O36").Value = 3700
$ws.Range("P36").Value = 3650
$ws.Range("Q36").Value = '$/kilo'
$ws.Range("R36").Value = 'Región del Maule'
$ws.Range("S36").Value = 3650
$ws.Range("T36").Value = 1

# Row 37
$ws.Range("D37").Value = 44979
$ws.Range("M37").Value = 100
$ws.Range("N37").Value = 2000
$ws.Range("O37").Value = 2200
$ws.Range("P37").Value = 2100
$ws.Range("S37").Value = 1050

# Row 38
$ws.Range("D38").Value = 44166
$ws.Range("M38").Value = 200
$ws.Range("N38").Value = 6000
$ws.Range("O38").Value = 6500
$ws.Range("P38").Value = 6250
$ws.Range("Q38").Value = '$/bandeja 12 canastillos 125 gramos'
$ws.Range("S38").Value = 4167
$ws.Range("T38").Value = 1.5

# Row 39
$ws.Range("D39").Value = 44176
$ws.Range("M39").Value = 300
$ws.Range("N39").Value = 5000
$ws.Range("O39").Value = 6000
$ws.Range("P39").Value = 5500
$ws.Range("S39").Value = 3667

# Row 40
$ws.Range("D40").Value = 44882
$ws.Range("M40").Value = 200
$ws.Range("N40").Value = 7500
$ws.Range("O40").Value = 8000
$ws.Range("P40").Value = 7750
$ws.Range("Q40").Value = '$/bandeja 12 canastillos 125 gramos'
$ws.Range("S40").Value = 5167
$ws.Range("T40").Value = 1.5

# Row 41
$ws.Range("D41").Value = 44936
$ws.Range("L41").Value = 'Primera'
$ws.Range("M41").Value = 400
$ws.Range("N41").Value = 2200
$ws.Range("O41").Value = 2500
$ws.Range("P41").Value = 2350
$ws.Range("Q41").Value = '$/bandeja 2 kilos'
$ws.Range("S41").Value = 1175
$ws.Range("T41").Value = 2

# Row 42
$ws.Range("D42").Value = 44922
$ws.Range("N42").Value = 3500
$ws.Range("O42").Value = 3800
$ws.Range("P42").Value = 3650
$ws.Range("Q42").Value = '$/bandeja 2 kilos'
$ws.Range("R42").Value = 'Provincia de Curicó'
$ws.Range("S42").Value = 1825
$ws.Range("T42").Value = 2

# Row 43
$ws.Range("D43").Value = 44960
$ws.Range("M43").Value = 400
$ws.Range("R43").Value = 'Provincia de Curicó'

# Row 44
$ws.Range("D44").Value = 44516
$ws.Range("M44").Value = 80
$ws.Range("N44").Value = 3700
$ws.Range("O44").Value = 3800
$ws.Range("P44").Value = 3750
$ws.Range("Q44").Value = '$/kilo'
$ws.Range("R44").Value = 'Región del Maule'
$ws.Range("S44").Value = 3750
$ws.Range("T44").Value = 1

# Row 45
$ws.Range("D45").Value = 44519
$ws.Range("M45").Value = 200
$ws.Range("N45").Value = 3700
$ws.Range("O45").Value = 3800
$ws.Range("P45").Value = 3750
$ws.Range("S45").Value = 3750

# Row 46
$ws.Range("D46").Value = 44995
$ws.Range("M46").Value = 200
$ws.Range("N46").Value = 2300
$ws.Range("O46").Value = 2500
$ws.Range("P46").Value = 2400
$ws.Range("Q46").Value = '$/bandeja 2 kilos'
$ws.Range("S46").Value = 1200
$ws.Range("T46").Value = 2

# Row 47
$ws.Range("D47").Value = 44904
$ws.Range("Q47").Value = '$/bandeja 2 kilos'
$ws.Range("R47").Value = 'Provincia de Curicó'
$ws.Range("S47").Value = 2625
$ws.Range("T47").Value = 2

# Row 48
$ws.Range("D48").Value = 44551
$ws.Range("M48").Value = 400
$ws.Range("N48").Value = 5000
$ws.Range("O48").Value = 5500
$ws.Range("P48").Value = 5250
$ws.Range("Q48").Value = '$/bandeja 12 canastillos 125 gramos'
$ws.Range("R48").Value = 'Región del Maule'
$ws.Range("S48").Value = 3500
$ws.Range("T48").Value = 1.5

# Row 49
$ws.Range("D49").Value = 44873
$ws.Range("N49").Value = 7500
$ws.Range("O49").Value = 8000
$ws.Range("P49").Value = 7750
$ws.Range("R49").Value = 'Provincia de Curicó'
$ws.Range("S49").Value = 5167

# Row 51
$ws.Range("D51").Value = 44901
$ws.Range("N51").Value = 5000
$ws.Range("O51").Value = 5500
$ws.Range("P51").Value = 5250
$ws.Range("S51").Value = 2625

# Row 52
$ws.Range("D52").Value = 44953
$ws.Range("M52").Value = 400
$ws.Range("N52").Value = 2000
$ws.Range("O52").Value = 2200
$ws.Range("P52").Value = 2100
$ws.Range("Q52").Value = '$/bandeja 2 kilos'
$ws.Range("S52").Value = 1050
$ws.Range("T52").Value = 2

# Row 53
$ws.Range("D53").Value = 44959

# Row 54
$ws.Range("D54").Value = 44911
$ws.Range("M54").Value = 600
$ws.Range("N54").Value = 4000
$ws.Range("O54").Value = 4200
$ws.Range("P54").Value = 4100
$ws.Range("Q54").Value = '$/bandeja 2 kilos'
$ws.Range("S54").Value = 2050
$ws.Range("T54").Value = 2

# Row 55
$ws.Range("D55").Value = 44970
$ws.Range("O55").Value = 2200
$ws.Range("P55").Value = 2100
$ws.Range("S55").Value = 1050

# Row 56
$ws.Range("D56").Value = 44544
$ws.Range("L56").Value = 'Primera'
$ws.Range("M56").Value = 400
$ws.Range("N56").Value = 5000
$ws.Range("O56").Value = 5500
$ws.Range("P56").Value = 5250
$ws.Range("Q56").Value = '$/bandeja 12 canastillos 125 gramos'
$ws.Range("R56").Value = 'Región del Maule'
$ws.Range("S56").Value = 3500
$ws.Range("T56").Value = 1.5

# Row 57
$ws.Range("D57").Value = 44950
$ws.Range("M57").Value = 200
$ws.Range("N57").Value = 2000
$ws.Range("O57").Value = 2200
$ws.Range("P57").Value = 2100
$ws.Range("S57").Value = 1050

# Row 58
$ws.Range("D58").Value = 44915
$ws.Range("M58").Value = 600
$ws.Range("N58").Value = 4000
$ws.Range("O58").Value = 4200
$ws.Range("P58").Value = 4100
$ws.Range("S58").Value = 2050

# Row 59
$ws.Range("D59").Value = 44914
$ws.Range("N59").Value = 3800
$ws.Range("O59").Value = 4000
$ws.Range("P59").Value = 3900
$ws.Range("S59").Value = 1950

# Row 60
$ws.Range("D60").Value = 44900
$ws.Range("M60").Value = 400
$ws.Range("N60").Value = 5500
$ws.Range("O60").Value = 6000
$ws.Range("P60").Value = 5750
$ws.Range("S60").Value = 2875

# Row 61
$ws.Range("D61").Value = 44880
$ws.Range("M61").Value = 300
$ws.Range("N61").Value = 7500
$ws.Range("O61").Value = 8000
$ws.Range("P61").Value = 7750
$ws.Range("Q61").Value = '$/bandeja 12 canastillos 125 gramos'
$ws.Range("S61").Value = 5167
$ws.Range("T61").Value = 1.5

# Row 62
$ws.Range("D62").Value = 44894
$ws.Range("M62").Value = 300
$ws.Range("N62").Value = 8000
$ws.Range("O62").Value = 8500
$ws.Range("P62").Value = 8250
$ws.Range("S62").Value = 4125

# Row 63
$ws.Range("D63").Value = 44957
$ws.Range("M63").Value = 200
$ws.Range("N63").Value = 2000
$ws.Range("O63").Value = 2200
$ws.Range("P63").Value = 2100
$ws.Range("S63").Value = 1050

# Row 64
$ws.Range("D64").Value = 44977
$ws.Range("M64").Value = 200
$ws.Range("N64").Value = 2000
$ws.Range("O64").Value = 2200
$ws.Range("P64").Value = 2100
$ws.Range("S64").Value = 1050

# Row 65
$ws.Range("D65").Value = 44921
$ws.Range("M65").Value = 200
$ws.Range("N65").Value = 3000
$ws.Range("O65").Value = 3500
$ws.Range("P65").Value = 3250
$ws.Range("S65").Value = 1625

# Row 66
$ws.Range("D66").Value = 44999
$ws.Range("M66").Value = 80
$ws.Range("N66").Value = 2300
$ws.Range("O66").Value = 2500
$ws.Range("P66").Value = 2400
$ws.Range("S66").Value = 1200

# Row 67
$ws.Range("D67").Value = 44939
$ws.Range("O67").Value = 2000
$ws.Range("P67").Value = 2000
$ws.Range("S67").Value = 1000

# Row 68
$ws.Range("D68").Value = 44939
$ws.Range("L68").Value = 'Segunda'
$ws.Range("M68").Value = 200
$ws.Range("N68").Value = 2200
$ws.Range("O68").Value = 2200
$ws.Range("P68").Value = 2200
$ws.Range("Q68").Value = '$/bandeja 2 kilos'
$ws.Range("S68").Value = 1100
$ws.Range("T68").Value = 2

# Row 69
$ws.Range("D69").Value = 44910
$ws.Range("N69").Value = 4000
$ws.Range("O69").Value = 4200
$ws.Range("P69").Value = 4100
$ws.Range("S69").Value = 2050
